$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.288.02'
$ws.Range("E2").Value = '  -1.15%  '
$ws.Range("D3").Value = '3.532.85'
$ws.Range("E3").Value = '  +0.22%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.35'
$ws.Range("E5").Value = '  -1.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.80'
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").Value = '3.534.13'
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.487'
$ws.Range("E9").Value = '  -0.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.125'
$ws.Range("E10").Value = '  +0.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.13'
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.380'
$ws.Range("E12").Value = '  -1.66%  '
$ws.Range("D13").Value = '4.141.69'
$ws.Range("E13").Value = '  +0.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.68'
$ws.Range("E14").Value = '  +0.00%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.119'
$ws.Range("E15").Value = '  +1.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000180'
$ws.Range("E16").Value = '  -1.10%  '
$ws.Range("D17").Value = '3.538.62'
$ws.Range("E17").Value = '  +0.28%  '
$ws.Range("D18").Value = '64.335.43'
$ws.Range("E18").Value = '  -1.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.82'
$ws.Range("E19").Value = '  -3.40%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.10'
$ws.Range("E20").Value = '  -2.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.64'
$ws.Range("E21").Value = '  -1.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '386.95'
$ws.Range("E22").Value = '  -1.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.577'
$ws.Range("E23").Value = '  -0.77%  '
$ws.Range("D24").Value = '3.679.12'
$ws.Range("E24").Value = '  +0.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.08'
$ws.Range("E25").Value = '  -1.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000116'
$ws.Range("E27").Value = '  +3.59%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.56'
$ws.Range("E28").Value = '  -1.71%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.52'
$ws.Range("E29").Value = '  -2.99%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.40'
$ws.Range("E31").Value = '  +0.86%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.24'
$ws.Range("E32").Value = '  -1.74%  '
$ws.Range("D33").Value = '3.545.96'
$ws.Range("E33").Value = '  +0.27%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.65'
$ws.Range("E35").Value = '  -1.80%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.146'
$ws.Range("E36").Value = '  +1.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.39'
$ws.Range("E37").Value = '  +1.41%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.56'
$ws.Range("E38").Value = '  -0.44%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.92'
$ws.Range("E39").Value = '  -0.54%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '160.79'
$ws.Range("E40").Value = '  -4.46%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0791'
$ws.Range("E41").Value = '  -2.55%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.817'
$ws.Range("E42").Value = '  -0.72%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.51'
$ws.Range("E43").Value = '  +2.72%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.11'
$ws.Range("E45").Value = '  -2.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.21'
$ws.Range("E46").Value = '  -4.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.42'
$ws.Range("E47").Value = '  -0.42%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.62'
$ws.Range("E48").Value = '  -2.33%  '
$ws.Range("D49").Value = '2.476.41'
$ws.Range("E49").Value = '  +2.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.84'
$ws.Range("E50").Value = '  -0.95%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.905'
$ws.Range("E51").Value = '  -0.42%  '
